$d = $word.ActiveDocument

$replacements = @(
    @{old="585÷4=146, 1"; new="745÷2=372, 1"},
    @{old="485÷4=121, 1"; new="293÷4=73, 1"},
    @{old="587÷6=97, 5"; new="288÷2=144, 0"},
    @{old="818÷6=136, 2"; new="308÷9=34, 2"},
    @{old="547÷8=68, 3"; new="892÷7=127, 3"},
    @{old="766÷8=95, 6"; new="877÷3=292, 1"},
    @{old="537÷8=67, 1"; new="965÷7=137, 6"},
    @{old="933÷4=233, 1"; new="154÷4=38, 2"},
    @{old="211÷6=35, 1"; new="360÷7=51, 3"},
    @{old="820÷5=164, 0"; new="262÷2=131, 0"},
    @{old="268÷8=33, 4"; new="762÷3=254, 0"},
    @{old="281÷9=31, 2"; new="771÷3=257, 0"},
    @{old="383÷5=76, 3"; new="282÷7=40, 2"},
    @{old="475÷6=79, 1"; new="399÷6=66, 3"},
    @{old="494÷8=61, 6"; new="238÷5=47, 3"},
    @{old="996÷6=166, 0"; new="685÷2=342, 1"},
    @{old="653÷6=108, 5"; new="477÷3=159, 0"},
    @{old="331÷2=165, 1"; new="689÷7=98, 3"},
    @{old="860÷7=122, 6"; new="509÷2=254, 1"},
    @{old="486÷7=69, 3"; new="890÷3=296, 2"},
    @{old="677÷4=169, 1"; new="358÷4=89, 2"},
    @{old="132÷6=22, 0"; new="154÷8=19, 2"},
    @{old="300÷4=75, 0"; new="803÷2=401, 1"},
    @{old="961÷9=106, 7"; new="262÷2=131, 0"},
    @{old="521÷9=57, 8"; new="454÷4=113, 2"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
